$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.828.90"
$ws.Range("E2").Value = "  +4.28%  "

$ws.Range("D3").Value = "1.875.71"
$ws.Range("E3").Value = "  +3.22%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "277.26"
$ws.Range("E5").Value = "  -0.11%  "

$ws.Range("D6").Value = "0.9999"

$ws.Range("D7").Value = "0.5279"
$ws.Range("E7").Value = "  +3.73%  "

$ws.Range("D8").Value = "0.3414"
$ws.Range("E8").Value = "  -3.21%  "

$ws.Range("D9").Value = "0.06942"
$ws.Range("E9").Value = "  +4.12%  "

$ws.Range("D10").Value = "20.07"
$ws.Range("E10").Value = "  +0.09%  "

$ws.Range("D11").Value = "0.8040"
$ws.Range("E11").Value = "  -2.89%  "

$ws.Range("D12").Value = "0.07718"
$ws.Range("E12").Value = "  -2.05%  "

$ws.Range("D13").Value = "1.880.45"
$ws.Range("E13").Value = "  +3.48%  "

$ws.Range("D14").Value = "5.185"
$ws.Range("E14").Value = "  +2.25%  "

$ws.Range("D15").Value = "90.26"
$ws.Range("E15").Value = "  +3.10%  "

$ws.Range("D16").Value = "14.58"
$ws.Range("E16").Value = "  +3.32%  "

$ws.Range("D17").Value = "0.9993"
$ws.Range("E17").Value = "  -0.05%  "

$ws.Range("D18").Value = "0.000008045"
$ws.Range("E18").Value = "  +0.22%  "

$ws.Range("D19").Value = "0.9999"

$ws.Range("D20").Value = "26.878.30"
$ws.Range("E20").Value = "  +4.27%  "

$ws.Range("D21").Value = "2.103.10"
$ws.Range("E21").Value = "  +2.68%  "

$ws.Range("D22").Value = "4.752"
$ws.Range("E22").Value = "  +0.24%  "

$ws.Range("D23").Value = "10.04"
$ws.Range("E23").Value = "  +0.35%  "

$ws.Range("D24").Value = "6.171"
$ws.Range("E24").Value = "  +1.28%  "

$ws.Range("D25").Value = "2.385"
$ws.Range("E25").Value = "  +8.38%  "

$ws.Range("D26").Value = "146.52"
$ws.Range("E26").Value = "  +2.98%  "

$ws.Range("D27").Value = "17.32"
$ws.Range("E27").Value = "  +1.23%  "

$ws.Range("D28").Value = "1.659"
$ws.Range("E28").Value = "  -0.96%  "

$ws.Range("D29").Value = "113.70"
$ws.Range("E29").Value = "  +3.83%  "

$ws.Range("E30").Value = "  +0.51%  "

$ws.Range("D31").Value = "4.311"
$ws.Range("E31").Value = "  +1.88%  "

$ws.Range("D32").Value = "0.08904"
$ws.Range("E32").Value = "  +1.46%  "

$ws.Range("D33").Value = "0.04927"
$ws.Range("E33").Value = "  +0.83%  "

$ws.Range("D34").Value = "1.174"
$ws.Range("E34").Value = "  +3.18%  "

$ws.Range("D35").Value = "0.7261"
$ws.Range("E35").Value = "  -0.21%  "

$ws.Range("D36").Value = "2.872"
$ws.Range("E36").Value = "  +0.13%  "

$ws.Range("D37").Value = "3.280"
$ws.Range("E37").Value = "  +5.00%  "

$ws.Range("D38").Value = "2.339"
$ws.Range("E38").Value = "  -1.40%  "

$ws.Range("D39").Value = "0.01858"
$ws.Range("E39").Value = "  +0.31%  "

$ws.Range("D40").Value = "0.5132"
$ws.Range("E40").Value = "  -0.86%  "

$ws.Range("D41").Value = "0.9546"
$ws.Range("E41").Value = "  -1.00%  "

$ws.Range("D42").Value = "116.01"
$ws.Range("E42").Value = "  +5.09%  "

$ws.Range("D43").Value = "6.150"
$ws.Range("E43").Value = "  -1.09%  "

$ws.Range("D44").Value = "8.113"
$ws.Range("E44").Value = "  +1.05%  "

$ws.Range("D45").Value = "0.9992"
$ws.Range("E45").Value = "  -0.10%  "

$ws.Range("D46").Value = "0.4474"
$ws.Range("E46").Value = "  -1.71%  "

$ws.Range("D47").Value = "0.1341"
$ws.Range("E47").Value = "  -1.76%  "

$ws.Range("D48").Value = "9.346"
$ws.Range("E48").Value = "  +0.96%  "

$ws.Range("D49").Value = "36.25"
$ws.Range("E49").Value = "  -0.69%  "

$ws.Range("D50").Value = "0.05936"
$ws.Range("E50").Value = "  +1.64%  "

$ws.Range("E51").Value = "  -0.75%  "
